# Insert a new weekly price record for "Pepino ensalada" (Vega Modelo de
# Temuco) at row 487. This pushes the existing data rows 487-524 down by
# one row (524 -> 525) and grows the used range from A1:R524 to A1:R525,
# matching the upstream weekly data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 487 (shifts 487..524 -> 488..525)
$ws.Rows(487).Insert()

# Populate the newly inserted row with the new record's values.
$ws.Range("A487").Value = 10
$ws.Range("B487").Value = "Vega Modelo de Temuco"
$ws.Range("C487").Value = "La Araucanía"
$ws.Range("D487").Value = 44746
$ws.Range("E487").Value = 9
$ws.Range("F487").Value = 100112043
$ws.Range("G487").Value = "Pepino ensalada"
$ws.Range("H487").Value = "Sin especificar"
$ws.Range("I487").Value = "Primera"
$ws.Range("J487").Value = 350
$ws.Range("K487").Value = 20000
$ws.Range("L487").Value = 20000
$ws.Range("M487").Value = 20000
$ws.Range("N487").Value = "$/caja 60 unidades"
$ws.Range("O487").Value = "Región de Arica y Parinacota"
$ws.Range("P487").Value = 333
$ws.Range("Q487").Value = 60
$ws.Range("R487").Value = "Hortaliza"
